$wb = $excel.ActiveWorkbook

# --- Add the new "TestCase2" worksheet (Login test case) ---
# With "TestCase1" as the active sheet, Worksheets.Add() inserts directly
# before it, landing exactly between "Test Cases" and "TestCase1".
$wsTC1 = $wb.Worksheets.Item("TestCase1")
$wsTC1.Activate()
$wsTC2 = $wb.Worksheets.Add()
$wsTC2.Name = "TestCase2"

# --- Update the "Test Cases" summary sheet with the new rows ---
# (shared-string order must match: TestCase2, Login are referenced first)
$wsSummary = $wb.Worksheets.Item("Test Cases")
$wsSummary.Range("A3").Value = "TestCase2"
$wsSummary.Range("B3").Value = "Login"
$wsSummary.Range("C3").Value = "y"
$wsSummary.Range("D3").Value = "PASS"

# --- Populate the new TestCase2 sheet ---
$wsTC2.Range("A1").Value = "Username"
$wsTC2.Range("B1").Value = "Password"
$wsTC2.Range("C1").Value = "Runmode"
$wsTC2.Range("D1").Value = "Results"
$wsTC2.Range("E1").Value = "Errors"

$wsTC2.Range("A2").Value = "Admin"
$wsTC2.Range("B2").Value = "admin123"
$wsTC2.Range("C2").Value = "y"
$wsTC2.Range("D2").Value = "PASS"

# --- Add the Logout row to the summary sheet ---
$wsSummary.Range("A4").Value = "TestCase3"
$wsSummary.Range("B4").Value = "Logout"
$wsSummary.Range("C4").Value = "y"
$wsSummary.Range("D4").Value = "PASS"

# --- Restore the originally active sheet / selection ---
$wsTC1.Activate()

Write-Host "done"
